$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (row 2-4, columns B-E) with new gain numbers
$ws.Range("B2").Value = 3.5362773999999999
$ws.Range("C2").Value = 10.028527
$ws.Range("D2").Value = 0.003
$ws.Range("E2").Value = 0.003

$ws.Range("B3").Value = 0.28603200000000001
$ws.Range("C3").Value = 0.118836
$ws.Range("D3").Value = 0.15
$ws.Range("E3").Value = 0.15

$ws.Range("B4").Value = 0.057206399999999998
$ws.Range("C4").Value = 0.023767
$ws.Range("D4").Value = 0.069000000000000006
$ws.Range("E4").Value = 0.069000000000000006

# Remove the centered style from column B (B1:B4) and B2:E4 block
$ws.Range("B1:B4").Style = "Normal"
$ws.Range("B2:E4").Style = "Normal"

# Set column B width (bestFit-like), matches new <col> entry
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# Update the active cell selection
$ws.Range("D8").Select()
